$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Rows 21 and 22 swap places (Uniswap <-> Avalanche)
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"

# Updated Price (D) and Volume(1h) (E) values
Set-TextValue $ws.Range("D2") "24.858.66"
$ws.Range("E2").Value = "  +2.02%  "
Set-TextValue $ws.Range("D3") "1.665.28"
$ws.Range("E3").Value = "  +1.24%  "
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  -0.25%  "
Set-TextValue $ws.Range("D5") "329.84"
$ws.Range("E5").Value = "  +7.97%  "
Set-TextValue $ws.Range("D6") "0.9983"
$ws.Range("E6").Value = "  -0.05%  "
Set-TextValue $ws.Range("D7") "0.3647"
$ws.Range("E7").Value = "  +0.58%  "
Set-TextValue $ws.Range("D8") "47.19"
$ws.Range("E8").Value = "  -0.11%  "
Set-TextValue $ws.Range("D9") "0.3244"
$ws.Range("E9").Value = "  -0.57%  "
Set-TextValue $ws.Range("D10") "1.142"
$ws.Range("E10").Value = "  +2.80%  "
Set-TextValue $ws.Range("D11") "0.07069"
$ws.Range("E11").Value = "  +2.25%  "
Set-TextValue $ws.Range("D12") "0.9990"
$ws.Range("E12").Value = "  -0.30%  "
Set-TextValue $ws.Range("D13") "6.074"
$ws.Range("E13").Value = "  +2.22%  "
Set-TextValue $ws.Range("D14") "19.66"
$ws.Range("E14").Value = "  +3.04%  "
Set-TextValue $ws.Range("D15") "1.664.88"
$ws.Range("E15").Value = "  +0.95%  "
Set-TextValue $ws.Range("D16") "6.607"
$ws.Range("E16").Value = "  +1.08%  "
Set-TextValue $ws.Range("D17") "0.00001054"
$ws.Range("E17").Value = "  +1.00%  "
Set-TextValue $ws.Range("D18") "0.06600"
$ws.Range("E18").Value = "  +1.68%  "
Set-TextValue $ws.Range("D19") "0.9989"
$ws.Range("E19").Value = "  -0.07%  "
Set-TextValue $ws.Range("D20") "78.82"
$ws.Range("E20").Value = "  +2.65%  "
Set-TextValue $ws.Range("D21") "15.87"
$ws.Range("E21").Value = "  +1.17%  "
Set-TextValue $ws.Range("D22") "5.936"
$ws.Range("E22").Value = "  +0.84%  "
Set-TextValue $ws.Range("D23") "12.51"
$ws.Range("E23").Value = "  +3.30%  "
Set-TextValue $ws.Range("D24") "24.857.04"
$ws.Range("E24").Value = "  +1.95%  "
Set-TextValue $ws.Range("D25") "2.444"
$ws.Range("E25").Value = "  +1.58%  "
Set-TextValue $ws.Range("D26") "2.421"
$ws.Range("E26").Value = "  +3.98%  "
Set-TextValue $ws.Range("D27") "148.81"
$ws.Range("E27").Value = "  +2.36%  "
Set-TextValue $ws.Range("D28") "18.69"
$ws.Range("E28").Value = "  +1.30%  "
Set-TextValue $ws.Range("D29") "1.846.54"
$ws.Range("E29").Value = "  +0.78%  "
Set-TextValue $ws.Range("D30") "125.62"
$ws.Range("E30").Value = "  +0.79%  "
Set-TextValue $ws.Range("D31") "1.185"
$ws.Range("E31").Value = "  +3.31%  "
Set-TextValue $ws.Range("D32") "4.062"
$ws.Range("E32").Value = "  +0.27%  "
Set-TextValue $ws.Range("D33") "5.741"
$ws.Range("E33").Value = "  +2.58%  "
Set-TextValue $ws.Range("D34") "0.08488"
$ws.Range("E34").Value = "  +2.05%  "
Set-TextValue $ws.Range("D35") "1.641"
$ws.Range("E35").Value = "  -1.96%  "
Set-TextValue $ws.Range("D36") "12.19"
$ws.Range("E36").Value = "  -0.56%  "
Set-TextValue $ws.Range("D37") "5.163"
$ws.Range("E37").Value = "  +0.03%  "
Set-TextValue $ws.Range("D38") "0.02262"
$ws.Range("E38").Value = "  +2.70%  "
Set-TextValue $ws.Range("D39") "0.06046"
$ws.Range("E39").Value = "  +0.32%  "
Set-TextValue $ws.Range("D40") "1.230"
$ws.Range("E40").Value = "  +2.35%  "
Set-TextValue $ws.Range("D41") "0.2085"
$ws.Range("E41").Value = "  +2.58%  "
Set-TextValue $ws.Range("D42") "8.231"
$ws.Range("E42").Value = "  +0.33%  "
Set-TextValue $ws.Range("D43") "0.9982"
$ws.Range("E43").Value = "  -0.15%  "
Set-TextValue $ws.Range("D44") "0.5933"
$ws.Range("E44").Value = "  +1.68%  "
Set-TextValue $ws.Range("D45") "13.48"
$ws.Range("E45").Value = "  +6.56%  "
$ws.Range("E46").Value = "  +3.48%  "
Set-TextValue $ws.Range("D47") "0.5665"
$ws.Range("E47").Value = "  +1.62%  "
Set-TextValue $ws.Range("D48") "125.58"
$ws.Range("E48").Value = "  +3.42%  "
Set-TextValue $ws.Range("D49") "1.951"
$ws.Range("E49").Value = "  +1.07%  "
Set-TextValue $ws.Range("D50") "0.06986"
$ws.Range("E50").Value = "  +1.39%  "
Set-TextValue $ws.Range("D51") "1.187"
$ws.Range("E51").Value = "  +4.02%  "
